# Add a new "Additional Information" column (D) to Tabelle1, with a couple
# of remarks, so sheets/test-procedures can carry extra notes (used later to
# compare sheets).

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Tabelle1")
$ws2 = $wb.Worksheets.Item("Sheet2")

# Header for the new column
$ws1.Range("D1").Value = "Additional Information"

# Row 2 (TP-01): nothing special to note
$ws1.Range("D2").Value = "-"

# Row 3 (TP-03): a remark, emphasised in bold
$ws1.Range("D3").Value = "… aber nur wenn sie leuchtet"
$ws1.Range("D3").Font.Bold = $true

# Give the new column a sensible width (engine rounds to whole pixels
# internally; 34.25 is the input that round-trips to the intended ~35.17)
$ws1.Columns.Item(4).ColumnWidth = 34.25

# Touch Sheet2 (keeps its selection in sync with the comparison work) then
# return focus to Tabelle1, which stays the active/selected sheet.
$ws2.Range("C22").Select() | Out-Null
$ws1.Activate() | Out-Null
$ws1.Range("D9").Select() | Out-Null
